# "more procs, not used anymore?"
#
# Adds two new proc-table entries to Sheet1:
#   - row 7: report_summary_testrun / ./perftools/report/report-run-dir.tcl
#   - row 9: a free-form note (column A only) explaining the proc isn't needed
# and moves the active selection down to A10, widening column A to fit the
# new (longer) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row 8 is intentionally left blank, matching the source sheet)
$ws.Range("A7").Value = "report_summary_testrun"
$ws.Range("B7").Value = "./perftools/report/report-run-dir.tcl"

$ws.Range("A9").Value = "[2017-04-03 11:03:30] Deze zou nu niet meer nodig moeten zijn, met stacktraces aanvulling, zie ndv::source_once.tcl"

# Column A needs to grow a bit to accommodate "report_summary_testrun"
$ws.Columns.Item(1).ColumnWidth = 20

# Selection moved on to the newly added note row
$ws.Range("A10").Select() | Out-Null
